# Fruta / hortaliza, semanal
# Swap the Fecha/Volumen/Precio values between rows 2-3 and rows 6-7
# (the "week" of data previously in rows 6-7 moves to rows 2-3 and vice versa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowValues($rowA, $rowB) {
    $cols = @("D", "M", "N", "O", "P", "S")
    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

Swap-RowValues 2 6
Swap-RowValues 3 7
